$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
  "Andrew Nembhard",
  "Quentin Grimes",
  "Klay Thompson",
  "Tobias Harris",
  "Zach LaVine",
  "Guerschon Yabusele",
  "Joel Embiid",
  "Kyrie Irving",
  "Shai Gilgeous-Alexander",
  "CJ McCollum",
  "Lauri Markkanen",
  "Jordan Poole",
  "Keyonte George",
  "Jalen Williams",
  "John Collins",
  "RJ Barrett",
  "Jimmy Butler"
)

$positions = @(
  "PG,SG",
  "SG,SF",
  "SG,SF",
  "SF,PF",
  "SG,SF",
  "PF,C",
  "C",
  "PG,SG",
  "PG,SG",
  "PG,SG",
  "SF,PF",
  "PG,SG",
  "PG,SG",
  "SG,SF,PF,C",
  "PF,C",
  "SG,SF,PF",
  "SF,PF"
)

$teams = @(
  "Indiana Pacers",
  "Dallas Mavericks",
  "Dallas Mavericks",
  "Detroit Pistons",
  "Chicago Bulls",
  "Philadelphia 76ers",
  "Philadelphia 76ers",
  "Dallas Mavericks",
  "Oklahoma City Thunder",
  "New Orleans Pelicans",
  "Utah Jazz",
  "Washington Wizards",
  "Utah Jazz",
  "Oklahoma City Thunder",
  "Utah Jazz",
  "Toronto Raptors",
  "Miami Heat"
)

for ($i = 0; $i -lt $names.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
}

for ($i = 0; $i -lt $positions.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $positions[$i]
}

for ($i = 0; $i -lt $teams.Count; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $teams[$i]
}
